$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tiles")

# Replace formulas with plain values in column A, and update column B values
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 3

$ws.Range("A2").Value = 15
$ws.Range("B2").Value = 2

$ws.Range("A3").Value = 25
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 35
$ws.Range("B4").Value = 2

# Update selection to match diff (activeCell B5, sqref B5)
$ws.Range("B5").Select()
